$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so Excel does not
# auto-convert numeric-looking strings (e.g. "218.02") into numbers.
$ws.Range("D2,E2").NumberFormat = "@"
$ws.Range("D2").Value = "26.326.85"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3,E3").NumberFormat = "@"
$ws.Range("D3").Value = "1.679.00"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5,E5").NumberFormat = "@"
$ws.Range("D5").Value = "218.02"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6,E6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5268"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("D7,E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8,E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2693"
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9,E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06463"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11,E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07509"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12,E12").NumberFormat = "@"
$ws.Range("D12").Value = "1.707.64"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13,E13").NumberFormat = "@"
$ws.Range("D13").Value = "4.517"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14,E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5787"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15,E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008514"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16,E16").NumberFormat = "@"
$ws.Range("D16").Value = "64.75"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17,E17").NumberFormat = "@"
$ws.Range("D17").Value = "26.335.34"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18,E18").NumberFormat = "@"
$ws.Range("D18").Value = "4.921"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21,E21").NumberFormat = "@"
$ws.Range("D21").Value = "189.76"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22,E22").NumberFormat = "@"
$ws.Range("D22").Value = "6.199"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24,E24").NumberFormat = "@"
$ws.Range("D24").Value = "144.85"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.08%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.48%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29,E29").NumberFormat = "@"
$ws.Range("D29").Value = "1.362"
$ws.Range("E29").Value = "  +4.36%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31,E31").NumberFormat = "@"
$ws.Range("D31").Value = "3.587"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35,E35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6212"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("D36,E36").NumberFormat = "@"
$ws.Range("D36").Value = "2.409"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37,E37").NumberFormat = "@"
$ws.Range("D37").Value = "2.735"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38,E38").NumberFormat = "@"
$ws.Range("D38").Value = "6.297"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39,E39").NumberFormat = "@"
$ws.Range("D39").Value = "1.114.61"
$ws.Range("E39").Value = "  +3.58%  "
$ws.Range("D40,E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01620"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41,E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8731"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43,E43").NumberFormat = "@"
$ws.Range("D43").Value = "100.46"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44,E44").NumberFormat = "@"
$ws.Range("D44").Value = "1.828.90"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45,E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000110"
$ws.Range("E45").Value = "  -3.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.85"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.87%  "
$ws.Range("D48,E48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49,E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05271"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50,E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4294"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51,E51").NumberFormat = "@"
$ws.Range("D51").Value = "6.071"
$ws.Range("E51").Value = "  +2.02%  "
